# Update column E ("SE") on "Sheet 1" from one constant-per-group value to
# a distinct per-row standard-error value for rows 2-82 (data rows beneath
# the TP/Site/Category/percent_cover/SE header in row 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 5).Value = 0.530967541957364
$ws.Cells.Item(4, 5).Value = 1.39422716504704
$ws.Cells.Item(5, 5).Value = 4.02989238728049
$ws.Cells.Item(6, 5).Value = 0.276730178760349
$ws.Cells.Item(7, 5).Value = 0.917601131611494
$ws.Cells.Item(8, 5).Value = 0.167112329324454
$ws.Cells.Item(9, 5).Value = 2.36215743656134
$ws.Cells.Item(10, 5).Value = 1.1850635047136
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(12, 5).Value = 0.0449077270887141
$ws.Cells.Item(13, 5).Value = 1.49647450473634
$ws.Cells.Item(14, 5).Value = 2.34561224277111
$ws.Cells.Item(15, 5).Value = 0.200700488427863
$ws.Cells.Item(16, 5).Value = 1.48831781985096
$ws.Cells.Item(17, 5).Value = 0.125254104140451
$ws.Cells.Item(18, 5).Value = 1.98609821754965
$ws.Cells.Item(19, 5).Value = 0.447788167666546
$ws.Cells.Item(20, 5).Value = 0.00869565217391304
$ws.Cells.Item(21, 5).Value = 0.04698239201842
$ws.Cells.Item(22, 5).Value = 0.858531863069872
$ws.Cells.Item(23, 5).Value = 1.92924831876499
$ws.Cells.Item(24, 5).Value = 0.176849147497534
$ws.Cells.Item(25, 5).Value = 1.77929946557363
$ws.Cells.Item(26, 5).Value = 0.0778828708630368
$ws.Cells.Item(27, 5).Value = 1.68706378023432
$ws.Cells.Item(28, 5).Value = 0.766127160396429
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(30, 5).Value = 0.121251972422417
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 5).Value = 4.8236026844199
$ws.Cells.Item(33, 5).Value = 0.335504493719715
$ws.Cells.Item(34, 5).Value = 1.92661488055065
$ws.Cells.Item(35, 5).Value = 0.247039615099626
$ws.Cells.Item(36, 5).Value = 4.62553339450419
$ws.Cells.Item(37, 5).Value = 1.22690917449051
$ws.Cells.Item(38, 5).Value = 0.00645161290322581
$ws.Cells.Item(39, 5).Value = 0.108513003655676
$ws.Cells.Item(40, 5).Value = 1.78388153580438
$ws.Cells.Item(41, 5).Value = 2.33747262257355
$ws.Cells.Item(42, 5).Value = 0.397477756741269
$ws.Cells.Item(43, 5).Value = 1.18166207370352
$ws.Cells.Item(44, 5).Value = 0.105222622130968
$ws.Cells.Item(45, 5).Value = 2.69576926170895
$ws.Cells.Item(46, 5).Value = 0.348096047489825
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 5).Value = 0.0678845782032709
$ws.Cells.Item(49, 5).Value = 1.12997115708868
$ws.Cells.Item(50, 5).Value = 1.71588657484018
$ws.Cells.Item(51, 5).Value = 0.115253786509998
$ws.Cells.Item(52, 5).Value = 1.27202055637692
$ws.Cells.Item(53, 5).Value = 0.109036671019418
$ws.Cells.Item(54, 5).Value = 1.95987178935455
$ws.Cells.Item(55, 5).Value = 0.532860502061678
$ws.Cells.Item(56, 5).Value = 0
$ws.Cells.Item(57, 5).Value = 1.08770118945923
$ws.Cells.Item(58, 5).Value = 1.14091263147231
$ws.Cells.Item(59, 5).Value = 4.08521475467774
$ws.Cells.Item(60, 5).Value = 0.314662112269589
$ws.Cells.Item(61, 5).Value = 1.50041082809383
$ws.Cells.Item(62, 5).Value = 0.292825782293005
$ws.Cells.Item(63, 5).Value = 2.80040230346007
$ws.Cells.Item(64, 5).Value = 0.827072535870849
$ws.Cells.Item(65, 5).Value = 0.00666666666666667
$ws.Cells.Item(66, 5).Value = 0.013198422611557
$ws.Cells.Item(67, 5).Value = 1.01498702342351
$ws.Cells.Item(68, 5).Value = 2.62208033716671
$ws.Cells.Item(69, 5).Value = 0.307109883003833
$ws.Cells.Item(70, 5).Value = 1.12021676819953
$ws.Cells.Item(71, 5).Value = 0.134753067726657
$ws.Cells.Item(72, 5).Value = 1.58939287908512
$ws.Cells.Item(73, 5).Value = 0.589257296156633
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 5).Value = 0.148516074749844
$ws.Cells.Item(76, 5).Value = 0.824422493512126
$ws.Cells.Item(77, 5).Value = 1.99764567432513
$ws.Cells.Item(78, 5).Value = 0.331095140219898
$ws.Cells.Item(79, 5).Value = 1.19091226108559
$ws.Cells.Item(80, 5).Value = 0.0820708360276816
$ws.Cells.Item(81, 5).Value = 1.52967968150779
$ws.Cells.Item(82, 5).Value = 0.839103283179233

Write-Output "Updated E2:E82 with per-row SE values"
